$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("08-10-2021", "09-10-2021", "10-10-2021", "11-10-2021", "12-10-2021", "13-10-2021")

$row = 282
foreach ($d in $dates) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $d
    $cell.Style = "Normal"
    $ws.Cells.Item($row, 2).Value = 696
    $ws.Cells.Item($row, 3).Value = 3962
    $ws.Cells.Item($row, 4).Value = 59
    $row++
}
